$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(82).Insert()

$ws.Cells.Item(82, 1).Value = 10
$ws.Cells.Item(82, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(82, 3).Value = "La Araucanía"
$ws.Cells.Item(82, 4).Value = 44665
$ws.Cells.Item(82, 5).Value = 9
$ws.Cells.Item(82, 6).Value = 100112005
$ws.Cells.Item(82, 7).Value = "Puerro"
$ws.Cells.Item(82, 8).Value = "Azul de Maquehue"
$ws.Cells.Item(82, 9).Value = "Primera"
$ws.Cells.Item(82, 10).Value = 30
$ws.Cells.Item(82, 11).Value = 12000
$ws.Cells.Item(82, 12).Value = 12000
$ws.Cells.Item(82, 13).Value = 12000
$ws.Cells.Item(82, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(82, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(82, 16).Value = 1000
$ws.Cells.Item(82, 17).Value = 12
$ws.Cells.Item(82, 18).Value = "Hortaliza"
